$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current "Total" row (row 15) so it becomes row 16,
# and the new data row becomes row 14.
$ws.Rows.Item(15).Insert()

# Fill in the new row 14 with the new timesheet entry, copying the
# formatting from the row above (row 13) so it matches the other entries.
$ws.Range("A14").Value = "Added controllers and views of the database"

$ws.Range("B13").Copy()
$ws.Range("B14").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("B14").Value = 43358

$ws.Range("C13").Copy()
$ws.Range("C14").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("C14").Value = 2

$excel.CutCopyMode = 0

# Update the Total formula to include the new row.
$ws.Range("C16").Formula = "=SUM(C2:C15)"

# Update the current selection to mirror the edited workbook state.
$ws.Range("C15").Select()
